# "separete file for b750" - update configlet values for the new B750 site
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Management Address / Loopback interface label
$ws.Range("B6").Value = "10.129.2.75"
$ws.Range("C6").Value = "Lo2"

# Managemet Mask / Uplink Addr 1 / Uplink Addr 2 - cleared out
$ws.Range("B7").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("B9").Value = ""

# Spine interfaces
$ws.Range("B10").Value = "E3/12"
$ws.Range("B11").Value = "E3/12"

# VRF
$ws.Range("B21").Value = "none"

# ASN
$ws.Range("B22").Value = 65012.1042

# D_VLAN_a / V_VLAN_a / SEC_VLAN_a / SP_VLAN_a addresses
$ws.Range("B23").Value = "10.129.75.1"
$ws.Range("B25").Value = "10.129.250.1"
$ws.Range("B27").Value = "10.129.252.1"
$ws.Range("B29").Value = "10.129.244.1"

# LoopbackID / LinkAddrID last octets
$ws.Range("B35").Value = 75
$ws.Range("B36").Value = 151

# SNMPLOC
$ws.Range("B38").Value = "B750 2CT2 1st Floor FanLab Rack Room"

# Move the active selection like the author left it
$ws.Range("B23").Select()
